$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows right before the existing row 308 ("2021-06-11" group),
# which pushes that group (and every group after it) down by three rows. The
# data that previously overflowed past row 431 now lands on new rows 432:434.
$ws.Rows("308:310").Insert()

# Populate the three newly-inserted blank rows with the new weekly record
# (date 2021-10-20, serial 44489) using the same constant columns shared by
# every row of this sheet (A, B, C, E, F, G, H, I, J, K).
$qualities = @('Especial', 'Primera', 'Segunda')
for ($i = 0; $i -lt 3; $i++) {
    $r = 308 + $i
    $ws.Cells.Item($r, 1).Value = 8
    $ws.Cells.Item($r, 2).Value = 'Terminal La Palmera de La Serena'
    $ws.Cells.Item($r, 3).Value = 'Coquimbo'
    $ws.Cells.Item($r, 4).Value = 44489
    $ws.Cells.Item($r, 5).Value = 4
    $ws.Cells.Item($r, 6).Value = 'Fruta'
    $ws.Cells.Item($r, 7).Value = 100108
    $ws.Cells.Item($r, 8).Value = 'Tropicales y subtropicales'
    $ws.Cells.Item($r, 9).Value = 100108002
    $ws.Cells.Item($r, 10).Value = 'Mango'
    $ws.Cells.Item($r, 11).Value = 'Sin especificar'
    $ws.Cells.Item($r, 12).Value = $qualities[$i]
    $ws.Cells.Item($r, 13).Value = 512
    $ws.Cells.Item($r, 14).Value = 6500
    $ws.Cells.Item($r, 15).Value = 7000
    $ws.Cells.Item($r, 16).Value = 6750
    $ws.Cells.Item($r, 17).Value = '$/bandeja 4 kilos'
    $ws.Cells.Item($r, 18).Value = 'Perú'
    $ws.Cells.Item($r, 19).Value = 1688
    $ws.Cells.Item($r, 20).Value = 4
}
